# Data van de loadtest: duplicate the "db / pp / train / predict" block
# (B3:E9) into a second block starting at column G (G3:J9), mirroring the
# same headers and values, then leave the selection where the user ended
# up after the paste.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing data table (headers row 3 + data rows 4-9, cols B-E)
# and paste it starting at G3 so it lands in G3:J9.
$ws.Range("B3:E9").Copy() | Out-Null
$ws.Range("G3").PasteSpecial() | Out-Null

# Copy/Paste in this runtime does not carry over the bold header style,
# so re-apply it explicitly to match the source header row (B3:E3 uses
# the bold style).
$ws.Range("G3:J3").Font.Bold = $true

$excel.CutCopyMode = $false

# Final selection left on the sheet after the edit.
$ws.Range("L13").Select() | Out-Null
